$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.674.36"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.517.82"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "317.29"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").Value = "95.69"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "36.35"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").Value = "2.904.99"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("E15").Value = "  +4.70%  "
$ws.Range("D16").Value = "2.533.73"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "42.688.22"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").Value = "0.0₃0975"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "71.50"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "253.16"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "3.00"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +12.92%  "
$ws.Range("D29").Value = "10.14"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "37.98"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").Value = "155.72"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  +4.61%  "
$ws.Range("D34").Value = "3.34"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").Value = "0.0787"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").Value = "24.15"
$ws.Range("E40").Value = "  -8.83%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "2.032.77"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").Value = "84.52"
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("D48").Value = "8.97"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "74.89"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "2.762.00"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("E51").Value = "  +0.36%  "
